# History & reports integrated
# - Remove the secondary AutoFilter criterion (Disease_Prediction = "Mange"),
#   keeping only the primary species = "Rabbit" filter. This unhides the
#   rows that were only hidden because they failed the second criterion.
# - Reset the sheet view: scroll back so the window is no longer pinned at
#   column V, and move the active selection to B268.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggling AutoFilter on field 22 (Disease_Prediction, the 22nd column of the
# filtered range A1:W433) removes that column's filter criteria while
# leaving the existing species filter (field 1) untouched.
$null = $ws.Range("A1:W433").AutoFilter(22)

# Move the selection to B268; this also resets the previously pinned
# top-left cell back to the default.
$null = $ws.Range("B268").Select()

Write-Host "AutoFilter criterion on Disease_Prediction removed; selection moved to B268"
